$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark. In the original document it sits
#    right after the inline picture (a few paragraphs below the bullet
#    list being edited); it is being relocated, not just deleted, so the
#    stale one has to go first.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Delete the whole bullet paragraph:
#    "Discover exactly what steps need to be done to generate a city and
#     generating different cities depending on some set of parameters."
#    (including its paragraph mark, so the remaining bullets close up).
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute( `
    "Discover exactly what steps need to be done to generate a city and generating different cities depending on some set of parameters.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.MoveEnd(1, 1)   # swallow the trailing paragraph mark too
    $rng.Delete()
}

# ---------------------------------------------------------------------
# 3. Change the trailing "?" to "." in the following bullet:
#    "Research what parameters and steps to use in the procedural
#     generation of cities?" -> "... generation of cities."
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("generation of cities?", $true, $false, $false, `
    $false, $false, $true, 1, $false, "generation of cities.", 2)

# ---------------------------------------------------------------------
# 4. Re-insert the "_GoBack" bookmark, collapsed, right at the end of
#    that same bullet paragraph (after the final run, before the
#    paragraph mark).
# ---------------------------------------------------------------------
$targetPara = $null
if ($found2) {
    $targetPara = $rng2.Paragraphs(1)
}

if ($targetPara -ne $null) {
    $paraEnd = $targetPara.Range.End

    # Collapsed ranges sitting exactly one position before a paragraph
    # mark are mis-handled by Bookmarks.Add, so a temporary marker
    # character is inserted first; that shifts the "end of paragraph
    # text" boundary away from the spot we actually want, lets the
    # bookmark be anchored correctly, and is then deleted again.
    $tailRng = $d.Range($paraEnd - 1, $paraEnd - 1)
    $tailRng.InsertAfter("X")

    $bmPos = $paraEnd - 1
    $bmRng = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRng)

    $markerRng = $d.Range($paraEnd - 1, $paraEnd)
    $markerRng.Delete()
}
